$d = $word.ActiveDocument

# --- 1) Remove the stray "_GoBack" bookmark from the empty paragraph
#        right after "HENINI Léo". It will be re-created later, attached
#        to the new final paragraph, further down in the document. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
}

# --- 2) Append the new paragraphs at the end of the document (after the
#        "D:\Wamp\bin\php\php5.6.38" paragraph, before the sectPr). Using
#        InsertXML lets us control the resulting OOXML precisely (style,
#        numbering, indentation, bookmark placement). ---
$endRange = $d.Content
$endRange.Collapse(0)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">' +
    '<pkg:xmlData>' +
      '<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' +
        '<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>' +
      '</Relationships>' +
    '</pkg:xmlData>' +
  '</pkg:part>' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +

          '<w:p>' +
            '<w:pPr>' +
              '<w:pStyle w:val="Paragraphedeliste"/>' +
              '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
            '</w:pPr>' +
            '<w:r><w:t xml:space="preserve">Il faut se connecter avec les comptes : </w:t></w:r>' +
          '</w:p>' +

          '<w:p>' +
            '<w:pPr>' +
              '<w:pStyle w:val="Paragraphedeliste"/>' +
              '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' +
            '</w:pPr>' +
            '<w:r><w:t>Login : admin  mdp : admin</w:t></w:r>' +
          '</w:p>' +

          '<w:p>' +
            '<w:pPr>' +
              '<w:pStyle w:val="Paragraphedeliste"/>' +
              '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' +
            '</w:pPr>' +
            '<w:r><w:t xml:space="preserve">Login : </w:t></w:r>' +
            '<w:r><w:t>p1812282</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve">  mdp : </w:t></w:r>' +
            '<w:r><w:t>hello</w:t></w:r>' +
          '</w:p>' +

          '<w:p>' +
            '<w:pPr>' +
              '<w:ind w:left="708"/>' +
            '</w:pPr>' +
            '<w:r><w:t>Les autres comptes n' + [char]0x2019 + 'ont pas de mots de passe hach' + [char]0x00E9 + 's, donc ne passe pas avec la fonction password_verify().</w:t></w:r>' +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
            '<w:bookmarkEnd w:id="0"/>' +
          '</w:p>' +

        '</w:body>' +
      '</w:document>' +
    '</pkg:xmlData>' +
  '</pkg:part>' +
'</pkg:package>'

$endRange.InsertXML($xml)
